$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous table content entirely (old table used A1:D8)
$ws.Range("A1:D8").Clear()

# Header row (row 1) - codenip/nom/prenom/parcours now start at column B
$ws.Range("B1").Value = "codenip"
$ws.Range("C1").Value = "nom"
$ws.Range("D1").Value = "prenom"
$ws.Range("E1").Value = "parcours"

# New "Competence" header columns (G1:K1)
$ws.Range("G1").Value = "Compétence 1"
$ws.Range("H1").Value = "Compétence 2"
$ws.Range("I1").Value = "Compétence 3"
$ws.Range("J1").Value = "Compétence 4"
$ws.Range("K1").Value = "Compétence 5"

# Row 2
$ws.Range("A2").Value = 12345
$ws.Range("B2").Value = "Etudiant 1"
$ws.Range("C2").Value = "Prenom1"
$ws.Range("D2").Value = "Parcours 1"
$ws.Range("G2").Value = "ADM"
$ws.Range("H2").Value = "ADM"

# Row 3
$ws.Range("A3").Value = 67890
$ws.Range("B3").Value = "Etudiant 2"
$ws.Range("C3").Value = "Prenom2"
$ws.Range("D3").Value = "Parcours 2"
$ws.Range("G3").Value = "ADM"

# Row 4
$ws.Range("A4").Value = 54321
$ws.Range("B4").Value = "Etudiant 3"
$ws.Range("C4").Value = "Prenom3"
$ws.Range("D4").Value = "Parcours 3"
$ws.Range("G4").Value = "NAR"

# Row 5
$ws.Range("A5").Value = 98765
$ws.Range("B5").Value = "Etudiant 4"
$ws.Range("C5").Value = "Prenom4"
$ws.Range("D5").Value = "Parcours 4"

# Row 6
$ws.Range("A6").Value = 13579
$ws.Range("B6").Value = "Etudiant 5"
$ws.Range("C6").Value = "Prenom5"
$ws.Range("D6").Value = "Parcours 5"

# The authored sheet's used range extends to column Q (A1:Q6) even though no
# values live past column K; touch Q6 with a no-op style so the sheet
# dimension/used-range grows to match without altering any shared styles.
$ws.Range("Q6").Font.Name = "Calibri"
